$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.130.86"
$ws.Range("E2").Value = '  +5.32%  '
$ws.Range("D3").Value = "'2.758.17"
$ws.Range("E3").Value = '  +4.52%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'581.36"
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = "'157.87"
$ws.Range("E6").Value = '  +9.40%  '
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = '  +4.73%  '
$ws.Range("D9").Value = "'2.757.38"
$ws.Range("E9").Value = '  +4.00%  '
$ws.Range("D10").Value = "'6.76"
$ws.Range("E10").Value = '  +2.84%  '
$ws.Range("D11").Value = "'0.112"
$ws.Range("E11").Value = '  +3.82%  '
$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = '  +4.63%  '
$ws.Range("D13").Value = "'0.159"
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").Value = "'3.224.03"
$ws.Range("E14").Value = '  +3.78%  '
$ws.Range("D15").Value = "'27.24"
$ws.Range("E15").Value = '  +3.81%  '
$ws.Range("D16").Value = "'64.008.44"
$ws.Range("E16").Value = '  +5.15%  '
$ws.Range("D18").Value = "'2.753.90"
$ws.Range("E18").Value = '  +3.97%  '
$ws.Range("D19").Value = "'12.06"
$ws.Range("E19").Value = '  +3.97%  '
$ws.Range("D20").Value = "'4.95"
$ws.Range("E20").Value = '  +4.71%  '
$ws.Range("D21").Value = "'363.73"
$ws.Range("E21").Value = '  +3.34%  '
$ws.Range("D22").Value = "'6.97"
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("D23").Value = "'0.549"
$ws.Range("E23").Value = '  +4.27%  '
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = "'67.13"
$ws.Range("E25").Value = '  +5.08%  '
$ws.Range("D26").Value = "'0.173"
$ws.Range("E26").Value = '  +6.53%  '
$ws.Range("D27").Value = "'8.64"
$ws.Range("E27").Value = '  +2.88%  '
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").Value = "'0.0₃0923"
$ws.Range("E29").Value = '  +14.09%  '
$ws.Range("D30").Value = "'2.02"
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("D31").Value = "'7.23"
$ws.Range("E31").Value = '  +7.15%  '
$ws.Range("D32").Value = "'1.27"
$ws.Range("E32").Value = '  +17.89%  '
$ws.Range("D33").Value = "'173.56"
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").Value = "'20.65"
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("D36").Value = "'4.94"
$ws.Range("E36").Value = '  +7.58%  '
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = '  +11.10%  '
$ws.Range("D38").Value = "'1.83"
$ws.Range("E38").Value = '  +8.38%  '
$ws.Range("D39").Value = "'1.02"
$ws.Range("E39").Value = '  +12.19%  '
$ws.Range("D40").Value = "'4.30"
$ws.Range("E40").Value = '  +4.45%  '
$ws.Range("D41").Value = "'337.41"
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("D42").Value = "'6.16"
$ws.Range("E42").Value = '  +16.06%  '
$ws.Range("D43").Value = "'39.47"
$ws.Range("E43").Value = '  +3.30%  '
$ws.Range("D44").Value = "'21.93"
$ws.Range("E44").Value = '  +8.06%  '
$ws.Range("D45").Value = "'22.25"
$ws.Range("E45").Value = '  +5.50%  '
$ws.Range("D46").Value = "'0.0604"
$ws.Range("E46").Value = '  +5.24%  '
$ws.Range("D47").Value = "'0.649"
$ws.Range("E47").Value = '  +3.83%  '
$ws.Range("D48").Value = "'0.0261"
$ws.Range("E48").Value = '  +4.15%  '
$ws.Range("D49").Value = "'137.32"
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").Value = "'0.103"
$ws.Range("E50").Value = '  +3.20%  '
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = '  +0.28%  '
